# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" and "Priority" columns for the rows
# whose Status is "Ready for handoff" on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 11, 13, 14)

# zh-cn sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZh.Range("H$r").Value = "2016-08-16 22:20:28"
    $wsZh.Range("E$r").Value = "ht"
}

# de-de sheet
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDe.Range("H$r").Value = "2016-08-16 22:20:34"
    $wsDe.Range("E$r").Value = "ht"
}

# Overview sheet mirrors the de-de "Latest Handoff Datetime" in its
# "Latest HO Xliff Generate Date" column (shared text in the workbook),
# so it must be kept in sync with the de-de update above.
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-16 22:20:34"
}
